$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row values (also establishes shared-string order) ---
$ws.Range("A1").Value = "Column1"
$ws.Range("B1").Value = "ColumnAfter1"
$ws.Range("C1").Value = "Column3"
$ws.Range("E1").Value = "Column4"

# --- Row 2 numeric / formatted cells (establishes numFmt style order) ---
$ws.Range("C2").Value = 3.5
$ws.Range("C2").NumberFormat = "0.00"
$ws.Range("D2").NumberFormat = "0.00"

$ws.Range("B2").Value = "2,5"
$ws.Range("B2").NumberFormat = "@"

$ws.Range("A2").Value = 1
$ws.Range("E2").Value = "four"

$ws.Range("F1").Value = "Column5Empty"
$ws.Range("G2").Value = "value6Filled"

$ws.Range("A3").Value = 10
$ws.Range("B3").Value = 20.5
$ws.Range("C3").Value = 30.5
$ws.Range("E3").Value = "fourX"
$ws.Range("G3").Value = "value6FilledX"

# --- Header row bold (establishes bold-font style last) ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true
$ws.Range("C1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true
$ws.Range("F1").Font.Bold = $true

# --- Column widths ---
# (COM ColumnWidth maps to stored OOXML width as `w + 5/6`, quantized to the
# nearest 1/6; these inputs land exactly/closest on the target widths of
# 12.5 and 21.1640625 respectively.)
$ws.Columns.Item(2).ColumnWidth = 11.6666666666667
$ws.Columns.Item(6).ColumnWidth = 20.3333333333333

# --- Selection ---
$ws.Range("C4").Select()
